$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New student rows (3-5) ---------------------------------------------
$ws.Range("A3").Value = "dhadhkasjdas"
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 34
$ws.Range("D3").Value = 5
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 3
$ws.Range("AQ3").Value = "X"

$ws.Range("A4").Value = "dfdsfsdf"
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 345
$ws.Range("E4").Value = 35
$ws.Range("F4").Value = 34
$ws.Range("G4").Value = 45345
$ws.Range("AQ4").Value = "X"

$ws.Range("A5").Value = "dsfsdfs"
$ws.Range("B5").Value = 1
$ws.Range("G5").Value = 345
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 5
$ws.Range("AN5").Value = 3

# --- Column G widened to fit the new 5-digit value (45345) -------------
$ws.Range("G1").ColumnWidth = 5.14

# --- Selection left on column C (whole column), as in the saved file ---
$ws.Columns("C").Select() | Out-Null
